# feat: add 2022-Q1 data
#
# - Inserts a new worksheet "2022-Q1" (holding the two fund rows for the
#   new quarter) positioned right before the "总计" summary sheet.
# - Updates the "总计" sheet with a new leading row summarizing 2022-Q1
#   (2 holdings, 0.06 亿元), shifting the previous rows down.

$wb = $excel.ActiveWorkbook

function Set-TextValue($rng, $val) {
    # Forces the cell to be stored as text (preserving things like leading
    # zeros in fund codes, or "0.32" rather than 0.32) while keeping the
    # cell's style reset to the sheet default - matching cells that carry
    # no explicit style index.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Step 1: insert a new worksheet "2022-Q1" right before "总计" ---
$totalSheet = $wb.Worksheets.Item("总计")
$q1Sheet = $wb.Worksheets.Add($totalSheet)    # inserted before $totalSheet
$q1Sheet.Name = "2022-Q1"

# Use "2021-Q4" as a layout/style template - headers and row formatting are
# identical, so copy the whole block (leaving the empty A1 corner alone)
# and then overwrite the fund-specific data cells.
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q4Sheet.Range("B1:H3").Copy($q1Sheet.Range("B1:H3"))
$q4Sheet.Range("A2:A3").Copy($q1Sheet.Range("A2:A3"))

# Row 2: 000927 - 博时大中华亚太精选股票(QDII) - 美元现汇
Set-TextValue $q1Sheet.Cells.Item(2, 2) "000927"
Set-TextValue $q1Sheet.Cells.Item(2, 3) "博时大中华亚太精选股票(QDII) - 美元现汇"
Set-TextValue $q1Sheet.Cells.Item(2, 4) "0.32"
Set-TextValue $q1Sheet.Cells.Item(2, 5) "92.94"
Set-TextValue $q1Sheet.Cells.Item(2, 6) "10.06"
Set-TextValue $q1Sheet.Cells.Item(2, 7) "0.0322"

# Row 3: 050015 - 博时大中华亚太精选股票(QDII) -人民币
Set-TextValue $q1Sheet.Cells.Item(3, 2) "050015"
Set-TextValue $q1Sheet.Cells.Item(3, 3) "博时大中华亚太精选股票(QDII) -人民币"
Set-TextValue $q1Sheet.Cells.Item(3, 4) "0.32"
Set-TextValue $q1Sheet.Cells.Item(3, 5) "92.94"
Set-TextValue $q1Sheet.Cells.Item(3, 6) "10.06"
Set-TextValue $q1Sheet.Cells.Item(3, 7) "0.0322"

# --- Step 2: update "总计" with a new leading row for 2022-Q1 ---
$totalSheet = $wb.Worksheets.Item("总计")

# Shift the existing two data rows down one row (bottom-up so we don't
# clobber data before it's been copied).
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A4:D4"))
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))

# Fix up the running row-index column (A) for the rows that moved down.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2

# Write the new first data row: 2022-Q1, 2 holdings, 0.06 亿元.
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 0.06
